# Outstandings.xlsx - "Add files via upload"
#
# Adds a new outstanding-purchase entry (Sr. No 3) to the
# "Purchase 22-23" sheet as row 30:
#   Sr.No=3, Date=04-Sep-2023 (45173), Invoice="PAN7404/23-24",
#   Vendor="Microciti", Bill amount=7198, Outstanding=E30
#
# Also mirrors the selection/view nudges that Excel records when a
# user scrolls/selects cells on both sheets while doing this edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# --- 1. Enter the new row's values -----------------------------------
$ws1.Range("A30").Value   = 3
$ws1.Range("B30").Value   = 45173
$ws1.Range("C30").Value   = "PAN7404/23-24"
$ws1.Range("D30").Value   = "Microciti"
$ws1.Range("E30").Value   = 7198
$ws1.Range("F30").Formula = "=E30"

# --- 2. Match the formatting used by the rest of the table -----------
# (row 2 is the first data row / start-of-group style used for A:F)
$ws1.Range("A2:F2").Copy()
$ws1.Range("A30:F30").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Update the on-screen selections that Excel persisted ---------
$ws2.Range("A28:G41").Select()

$ws1.Activate()
$ws1.Range("A31").Select()
